$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 308, shifting existing rows 308..372 down to 309..373.
$ws.Rows.Item(308).Insert()

# Populate the newly inserted row 308 with its data (same "fixed" attributes as its
# neighbours, new date/variety/volume/price figures).
$ws.Cells.Item(308, 1).Value = 1
$ws.Cells.Item(308, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(308, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(308, 4).Value = 44995
$ws.Cells.Item(308, 5).Value = 15
$ws.Cells.Item(308, 6).Value = "Fruta"
$ws.Cells.Item(308, 7).Value = 100108
$ws.Cells.Item(308, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(308, 9).Value = 100108006
$ws.Cells.Item(308, 10).Value = "Plátano"
$ws.Cells.Item(308, 11).Value = "Sin especificar"
$ws.Cells.Item(308, 12).Value = "Verde"
$ws.Cells.Item(308, 13).Value = 250
$ws.Cells.Item(308, 14).Value = 24000
$ws.Cells.Item(308, 15).Value = 25000
$ws.Cells.Item(308, 16).Value = 24600
$ws.Cells.Item(308, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(308, 18).Value = "Ecuador"
$ws.Cells.Item(308, 19).Value = 1230
$ws.Cells.Item(308, 20).Value = 20
